# Add the predefined term table
#
# The sheet previously held a couple of GeoMaterialConfidence example rows
# (row 2 "GeoMaterialConfidence/Heading2/Low/Prepop..." and row 3
# "DescriptionOfMapUnits/GeoMaterialConfidence/Low/2 Prepop 2..."). This
# replaces that sample data with the real predefined-term row
# (DescriptionOfMapUnits / ParagraphStyle / Heading / Heading Definition)
# and renames the sheet to match the workbook's purpose.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "PredefinedTerms"

# Row 1 (headers: DatasetName, FieldName, Term, Definition) is unchanged.

# Row 2 becomes the "Heading" predefined term entry.
$ws.Range("A2").Value = "DescriptionOfMapUnits"
$ws.Range("B2").Value = "ParagraphStyle"
$ws.Range("C2").Value = "Heading"
$ws.Range("D2").Value = "Heading Definition"

# Row 3 is removed entirely, shrinking the used range back to A1:D2.
$ws.Range("A3:D3").ClearContents()

# Leave the selection on the last edited cell, like the live edit session did.
$ws.Range("D2").Select()
